try {
  Write-Output ("const: " + $ppDateTimeFigureOut)
} catch { Write-Output ("err1: " + $_) }
try {
  Write-Output ("const2: " + [int]$ppDateTimeFigureOut)
} catch { Write-Output ("err2: " + $_) }
